$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = "Crawling in Python(request, BeautifulSoup, Selenium) (1)"
$ws.Range("E6").Value = "https://leedakyeong.tistory.com/entry/Crawling-in-Pythonrequest-BeautifulSoup-Selenium"

$ws.Range("D9").Value = "[해외DS] Microsoft, ‘이상한 행동’ 한 Bing AI 챗봇 제한한다"
$ws.Range("E9").Value = "https://pdsi.pabii.com/microsoft-limits-bing-ai-chat-generations-after-weird-behavior/#utm_source=rss&utm_medium=rss&utm_campaign=microsoft-limits-bing-ai-chat-generations-after-weird-behavior"

$ws.Range("D44").Value = "2023 MWC 망 사용료 논의"
$ws.Range("E44").Value = "https://engineering-ladder.tistory.com/120"

$ws.Range("D51").Value = "PyQt6와 PySide6의 라이선스 차이"
$ws.Range("E51").Value = "https://bskyvision.com/entry/PyQt6%EC%99%80-PySide6%EC%9D%98-%EB%9D%BC%EC%9D%B4%EC%84%A0%EC%8A%A4-%EC%B0%A8%EC%9D%B4"
